$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at the top of the data (rows 2 and 3), pushing the
# existing data rows (2-10) down to rows 4-12.
$ws.Rows("2:3").Insert(-4121, 0)

# The inserted rows pick up formatting (bold/centered) copied from row 1.
# Reset them back to the default "Normal" style used by the other data rows.
$ws.Range("A2:Q3").ClearFormats()

# Helper to write a text value into a cell while preventing Excel from
# auto-converting strings that look like dates (e.g. "2024-05-02") into
# date serial numbers, and without leaving a lingering custom style on
# the cell itself.
function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value2 = $text
    $rng.Style = "Normal"
}

# New row 2: 디앤디파마텍
Set-TextValue "A2" "2024-05-02"
$ws.Range("B2").Value2 = "디앤디파마텍"
$ws.Range("C2").Value2 = "코스닥"
$ws.Range("D2").Value2 = 363
$ws.Range("E2").Value2 = "한국"
$ws.Range("F2").Value2 = 363
$ws.Range("G2").Value2 = "-"
$ws.Range("H2").Value2 = "-"
$ws.Range("I2").Value2 = "-"
$ws.Range("J2").Value2 = "-"
$ws.Range("K2").Value2 = "대표"
$ws.Range("L2").Value2 = "-"
$ws.Range("M2").Value2 = 33000
$ws.Range("N2").Value2 = 100
Set-TextValue "O2" "2024-04-22"
Set-TextValue "P2" "2024-04-25"
$ws.Range("Q2").Value2 = 805400

# New row 3: 유안타제16호스팩
Set-TextValue "A3" "2024-05-02"
$ws.Range("B3").Value2 = "유안타제16호스팩"
$ws.Range("C3").Value2 = "코스닥"
$ws.Range("D3").Value2 = 103
$ws.Range("E3").Value2 = "유안타"
$ws.Range("F3").Value2 = 103
$ws.Range("G3").Value2 = "-"
$ws.Range("H3").Value2 = "-"
$ws.Range("I3").Value2 = "-"
$ws.Range("J3").Value2 = "-"
$ws.Range("K3").Value2 = "대표"
$ws.Range("L3").Value2 = "-"
$ws.Range("M3").Value2 = 2000
$ws.Range("N3").Value2 = 100
Set-TextValue "O3" "2024-04-22"
Set-TextValue "P3" "2024-04-25"
$ws.Range("Q3").Value2 = 3862500

Write-Output "done"
